$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Models")

# Fill columns E:G first (row by row), matching the order the new
# shared-string entries were introduced in, then fill column D.
$ws.Range("E1").Value = "shares_roc_var"
$ws.Range("F1").Value = "tech_var"
$ws.Range("G1").Value = "histend_var"

$ws.Range("E2").Value = "MSRC"
$ws.Range("F2").Value = "T2TI"
$ws.Range("G2").Value = "MEWG"

$ws.Range("E3").Value = "na"
$ws.Range("F3").Value = "HTTI"
$ws.Range("G3").Value = "HEWF"

$ws.Range("E4").Value = "na"
$ws.Range("F4").Value = "VTTI"
$ws.Range("G4").Value = "TEWS"

$ws.Range("E5").Value = "na"
$ws.Range("F5").Value = "na"
$ws.Range("G5").Value = "na"

$ws.Range("E6").Value = "na"
$ws.Range("F6").Value = "na"
$ws.Range("G6").Value = "na"

$ws.Range("E7").Value = "na"
$ws.Range("F7").Value = "na"
$ws.Range("G7").Value = "na"

$ws.Range("E8").Value = "na"
$ws.Range("F8").Value = "na"
$ws.Range("G8").Value = "na"

$ws.Range("E9").Value = "na"
$ws.Range("F9").Value = "na"
$ws.Range("G9").Value = "na"

$ws.Range("E10").Value = "na"
$ws.Range("F10").Value = "na"
$ws.Range("G10").Value = "na"

$ws.Range("D1").Value = "shares_var"
$ws.Range("D2").Value = "MEWS"
$ws.Range("D3").Value = "HEWS"
$ws.Range("D4").Value = "TEWS"
$ws.Range("D5").Value = "na"
$ws.Range("D6").Value = "IWS1"
$ws.Range("D7").Value = "IWS2"
$ws.Range("D8").Value = "IWS3"
$ws.Range("D9").Value = "IWS4"
$ws.Range("D10").Value = "IWS5"

# Move the selection/active tab to the Models sheet (was on C6TI before).
$ws.Range("E15").Select()
$ws.Activate()
